$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the whole bullet paragraph about the "th" background color fix
#    ("En la sección de grados ... #404040."). Its content was effectively
#    replaced by the (pre-existing) "Uso de etiqueta obsoleta center..."
#    paragraph shifting up, so we simply delete this paragraph outright.
# ---------------------------------------------------------------------------
$rngGrados = $d.Content
$found1 = $rngGrados.Find.Execute("En la sección de grados un contraste de colores con dificultades para la visión entre las etiquetas th, por ello he cambiado el color de fondo de los th a ")
if ($found1) {
    $rngGrados.Expand(4) | Out-Null   # wdParagraph -> whole paragraph incl. mark
    $rngGrados.Delete()
}

# ---------------------------------------------------------------------------
# 2) Fix the "Faltas de etiquetas formulario..." paragraph text: replace the
#    "...había que encerrar los input en etiquetas label." ending with
#    "...había que utilizar etiquetas label para referenciar  los input."
# ---------------------------------------------------------------------------
$rngForm = $d.Content
$found2 = $rngForm.Find.Execute("Faltas de etiquetas formulario en la página Solicitud.html")
if ($found2) {
    $rngForm.Expand(4) | Out-Null     # whole paragraph incl. mark
    $pStart = $rngForm.Start

    $run1Text = "Faltas de etiquetas formulario en la página Solicitud.html"
    $oldRun2Text = " y también en esta página había que encerrar los input en etiquetas label"
    $oldRun3Text = "."

    $run2Start = $pStart + $run1Text.Length
    $run2End = $run2Start + $oldRun2Text.Length

    # Replace the middle run's text in place.
    $rngRun2 = $d.Range($run2Start, $run2End)
    $rngRun2.Text = " y también en esta página había que "
    $afterFirst = $rngRun2.End

    # Insert the new third run right after it.
    $insertPt = $d.Range($afterFirst, $afterFirst)
    $insertPt.InsertAfter("utilizar etiquetas label para referenciar ")
    $afterSecond = $insertPt.End

    # The old trailing "." run now sits right after; turn it into the new
    # closing run text.
    $rngRun3 = $d.Range($afterSecond, $afterSecond + $oldRun3Text.Length)
    $rngRun3.Text = " los input."
}
